# "User stories #36 - #44, #51 - #57 + new requirments"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlCenter = -4108
$xlPasteFormats = -4122

# ---- New requirement text (rows 13-15, column B) ----
# Shared-string insertion order needs to be: Hebrew, New, Detailed logs, Data never deleted
$ws.Range("B13").Value = "User interfaces will be represented in Hebrew"
$ws.Range("G2").Value = "New"
$ws.Range("B14").Value = "Detailed logs will be written for each component"
$ws.Range("B15").Value = "Data will never be permanently deleted from the Database"

# ---- Priority (C) / Complexity (D) values for rows 2-15 ----
$priority   = @{2=5; 3=4; 4=4; 5=4; 6=4; 7=3; 8=5; 9=2; 10=3; 11=5; 12=5; 13=5; 14=5; 15=4}
$complexity = @{2=4; 3=2; 4=2; 5=3; 6=3; 7=3; 8=3; 9=4; 10=3; 11=3; 12=5; 13=1; 14=3; 15=1}

foreach ($r in 2..15) {
    $ws.Range("C$r").Value = $priority[$r]
    $ws.Range("D$r").Value = $complexity[$r]
    $ws.Range("G$r").Value = "New"
}

# ---- Apply the plain center/center alignment style to C,D,F,G (rows 2-17) and E (rows 16-17) ----
# (union of areas, applied per-area so the engine reuses a single style record)
$plainAreasRange = $excel.Union($ws.Range("C2:D17"), $ws.Range("F2:G17"), $ws.Range("E16:E17"))
foreach ($area in $plainAreasRange.Areas) {
    $area.HorizontalAlignment = $xlCenter
    $area.VerticalAlignment = $xlCenter
}

# ---- Build the date-formatted style once on E2, then propagate via copy/paste-special ----
# (setting NumberFormat to a date format on a multi-cell range creates a distinct style
#  per cell in this engine, so build it on a single cell and fan it out with PasteSpecial)
$ws.Range("E2").HorizontalAlignment = $xlCenter
$ws.Range("E2").VerticalAlignment = $xlCenter
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Value = 42690

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3:E15").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

foreach ($r in 3..15) {
    $ws.Range("E$r").Value = 42690
}

# E16/E17 stay blank but use the plain center style (already applied above via the union)

# ---- Update the selection shown in the sheet view ----
$ws.Range("B15").Select()
